$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Il16"
$ws.Range("C2").Value = "Grin2d"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.933820999999999
$ws.Range("H2").Value = 23.801463
$ws.Range("I2").Value = 0.3933990084177585
$ws.Range("J2").Value = 0.3933990084177585
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1653213333333333
$ws.Range("N2").Value = 0.495964
$ws.Range("O2").Value = 0.07076294988206985
$ws.Range("P2").Value = 0.07076294988206984
$ws.Range("Q2").Value = 1.311629866148
$ws.Range("R2").Value = 11.804668795332
$ws.Range("S2").Value = 0.02783807431632182
$ws.Range("T2").Value = 0.02783807431632181

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Il16"
$ws.Range("C3").Value = "Grin2d"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.933820999999999
$ws.Range("H3").Value = 23.801463
$ws.Range("I3").Value = 0.3933990084177585
$ws.Range("J3").Value = 0.3933990084177585
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.792633666666666
$ws.Range("N3").Value = 5.377901
$ws.Range("O3").Value = 0.7673059716707931
$ws.Range("P3").Value = 0.7673059716707932
$ws.Range("Q3").Value = 14.222434629907
$ws.Range("R3").Value = 128.001911669163
$ws.Range("S3").Value = 0.3018574084083147
$ws.Range("T3").Value = 0.3018574084083147

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Il16"
$ws.Range("C4").Value = "Grin2d"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.933820999999999
$ws.Range("H4").Value = 23.801463
$ws.Range("I4").Value = 0.3933990084177585
$ws.Range("J4").Value = 0.3933990084177585
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.3783146666666666
$ws.Range("N4").Value = 1.134944
$ws.Range("O4").Value = 0.161931078447137
$ws.Range("P4").Value = 0.161931078447137
$ws.Range("Q4").Value = 3.001480847007999
$ws.Range("R4").Value = 27.013327623072
$ws.Range("S4").Value = 0.06370352569312197
$ws.Range("T4").Value = 0.06370352569312197

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Il16"
$ws.Range("C5").Value = "Grin2d"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.58194133333333
$ws.Range("H5").Value = 34.745824
$ws.Range("I5").Value = 0.5742912823576415
$ws.Range("J5").Value = 0.5742912823576415
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1653213333333333
$ws.Range("N5").Value = 0.495964
$ws.Range("O5").Value = 0.07076294988206985
$ws.Range("P5").Value = 0.07076294988206984
$ws.Range("Q5").Value = 1.914741983815111
$ws.Range("R5").Value = 17.232677854336
$ws.Range("S5").Value = 0.04063854523118341
$ws.Range("T5").Value = 0.0406385452311834

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Il16"
$ws.Range("C6").Value = "Grin2d"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.58194133333333
$ws.Range("H6").Value = 34.745824
$ws.Range("I6").Value = 0.5742912823576415
$ws.Range("J6").Value = 0.5742912823576415
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.792633666666666
$ws.Range("N6").Value = 5.377901
$ws.Range("O6").Value = 0.7673059716707931
$ws.Range("P6").Value = 0.7673059716707932
$ws.Range("Q6").Value = 20.76217795949155
$ws.Range("R6").Value = 186.859601635424
$ws.Range("S6").Value = 0.4406571304314959
$ws.Range("T6").Value = 0.440657130431496

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Il16"
$ws.Range("C7").Value = "Grin2d"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.58194133333333
$ws.Range("H7").Value = 34.745824
$ws.Range("I7").Value = 0.5742912823576415
$ws.Range("J7").Value = 0.5742912823576415
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.3783146666666666
$ws.Range("N7").Value = 1.134944
$ws.Range("O7").Value = 0.161931078447137
$ws.Range("P7").Value = 0.161931078447137
$ws.Range("Q7").Value = 4.381618274872888
$ws.Range("R7").Value = 39.434564473856
$ws.Range("S7").Value = 0.09299560669496217
$ws.Range("T7").Value = 0.09299560669496217

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Il16"
$ws.Range("C8").Value = "Grin2d"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.6516016666666666
$ws.Range("H8").Value = 1.954805
$ws.Range("I8").Value = 0.03230970922460003
$ws.Range("J8").Value = 0.03230970922460003
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.1653213333333333
$ws.Range("N8").Value = 0.495964
$ws.Range("O8").Value = 0.07076294988206985
$ws.Range("P8").Value = 0.07076294988206984
$ws.Range("Q8").Value = 0.1077236563355556
$ws.Range("R8").Value = 0.96951290702
$ws.Range("S8").Value = 0.002286330334564622
$ws.Range("T8").Value = 0.002286330334564622

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Il16"
$ws.Range("C9").Value = "Grin2d"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.6516016666666666
$ws.Range("H9").Value = 1.954805
$ws.Range("I9").Value = 0.03230970922460003
$ws.Range("J9").Value = 0.03230970922460003
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.792633666666666
$ws.Range("N9").Value = 5.377901
$ws.Range("O9").Value = 0.7673059716707931
$ws.Range("P9").Value = 0.7673059716707932
$ws.Range("Q9").Value = 1.168083084922778
$ws.Range("R9").Value = 10.512747764305
$ws.Range("S9").Value = 0.02479143283098251
$ws.Range("T9").Value = 0.02479143283098252

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Il16"
$ws.Range("C10").Value = "Grin2d"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.6516016666666666
$ws.Range("H10").Value = 1.954805
$ws.Range("I10").Value = 0.03230970922460003
$ws.Range("J10").Value = 0.03230970922460003
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.3783146666666666
$ws.Range("N10").Value = 1.134944
$ws.Range("O10").Value = 0.161931078447137
$ws.Range("P10").Value = 0.161931078447137
$ws.Range("Q10").Value = 0.2465104673244444
$ws.Range("R10").Value = 2.21859420592
$ws.Range("S10").Value = 0.005231946059052895
$ws.Range("T10").Value = 0.005231946059052895

